$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.855.13'
$ws.Range("E2").Value = '  -0.59%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.498.67'
$ws.Range("E3").Value = '  -1.95%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.89'
$ws.Range("E5").Value = '  -1.11%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '197.31'
$ws.Range("E6").Value = '  +6.02%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.624'
$ws.Range("E7").Value = '  +1.08%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("E9").Value = '  -2.84%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.653'
$ws.Range("E10").Value = '  +1.43%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.15'
$ws.Range("E11").Value = '  +0.54%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000302'
$ws.Range("E12").Value = '  -2.50%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.58'
$ws.Range("E13").Value = '  +0.90%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.057.71'
$ws.Range("E14").Value = '  -1.83%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '602.79'
$ws.Range("E15").Value = '  +3.72%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.967.27'
$ws.Range("E16").Value = '  -0.54%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.07'
$ws.Range("E17").Value = '  +0.51%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.60'
$ws.Range("E18").Value = '  -0.76%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.495.14'
$ws.Range("E19").Value = '  -2.20%  '

$ws.Range("E20").Value = '  +0.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.988'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '18.36'
$ws.Range("E22").Value = '  +5.94%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '104.44'
$ws.Range("E23").Value = '  +10.02%  '

$ws.Range("E24").Value = '  -3.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.02'
$ws.Range("E25").Value = '  +3.47%  '

$ws.Range("E26").Value = '  +5.10%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.97'
$ws.Range("E27").Value = '  +0.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.74'
$ws.Range("E28").Value = '  +3.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.59'
$ws.Range("E29").Value = '  +3.94%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.49'
$ws.Range("E30").Value = '  +21.94%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.23'
$ws.Range("E31").Value = '  +2.73%  '

$ws.Range("E32").Value = '  +3.83%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.115'
$ws.Range("E33").Value = '  +1.04%  '

$ws.Range("E34").Value = '  -0.02%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.730.35'
$ws.Range("E35").Value = '  +5.49%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0₃0806'
$ws.Range("E36").Value = '  +2.73%  '

$ws.Range("E37").Value = '  -0.24%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '510.82'
$ws.Range("E38").Value = '  -3.95%  '

$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.99'
$ws.Range("E39").Value = '  -7.81%  '

$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.391'
$ws.Range("E40").Value = '  -3.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.61'
$ws.Range("E41").Value = '  -1.36%  '

$ws.Range("E42").Value = '  -0.47%  '

$ws.Range("E43").Value = '  +0.83%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0457'
$ws.Range("E44").Value = '  -0.60%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.83'
$ws.Range("E45").Value = '  -3.42%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.140'
$ws.Range("E46").Value = '  -0.46%  '

$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.32'
$ws.Range("E47").Value = '  -4.09%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.01'
$ws.Range("E48").Value = '  +0.36%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.74'
$ws.Range("E49").Value = '  -5.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.91'
$ws.Range("E50").Value = '  -3.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000241'
$ws.Range("E51").Value = '  -2.00%  '
